# Mortality_tables/Table_Summary.xlsx edit: add IfoA 92 Series rows (46-61)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phase 1: fill Datasource (A) and Table (B) columns for all new rows, row by row
$ws.Range("A46").Value = "IfoA 92 Series"
$ws.Range("B46").Value = "AF92"
$ws.Range("A47").Value = "IfoA 92 Series"
$ws.Range("B47").Value = "AM92"
$ws.Range("A48").Value = "IfoA 92 Series"
$ws.Range("B48").Value = "IFA92"
$ws.Range("A49").Value = "IfoA 92 Series"
$ws.Range("B49").Value = "IFL92"
$ws.Range("A50").Value = "IfoA 92 Series"
$ws.Range("B50").Value = "IMA92"
$ws.Range("A51").Value = "IfoA 92 Series"
$ws.Range("B51").Value = "IML92"
$ws.Range("A52").Value = "IfoA 92 Series"
$ws.Range("B52").Value = "PFA92"
$ws.Range("A53").Value = "IfoA 92 Series"
$ws.Range("B53").Value = "PFL92"
$ws.Range("A54").Value = "IfoA 92 Series"
$ws.Range("B54").Value = "PMA92"
$ws.Range("A55").Value = "IfoA 92 Series"
$ws.Range("B55").Value = "PML92"
$ws.Range("A56").Value = "IfoA 92 Series"
$ws.Range("B56").Value = "RFV92"
$ws.Range("A57").Value = "IfoA 92 Series"
$ws.Range("B57").Value = "RMV92"
$ws.Range("A58").Value = "IfoA 92 Series"
$ws.Range("B58").Value = "TF92"
$ws.Range("A59").Value = "IfoA 92 Series"
$ws.Range("B59").Value = "TM92"
$ws.Range("A60").Value = "IfoA 92 Series"
$ws.Range("B60").Value = "WA92"
$ws.Range("A61").Value = "IfoA 92 Series"
$ws.Range("B61").Value = "WL92"

# Phase 2: fill Table Description (C) and Datasource Location (E) columns, row by row
$ws.Range("C46").Value = "AF92: Permanent Assurances, females, combined"
$ws.Range("E46").Value = "Mortality_tables/92series.xls"
$ws.Range("C47").Value = "AM92: Permanent Assurances, males, combined"
$ws.Range("E47").Value = "Mortality_tables/92series.xls"
$ws.Range("C48").Value = "IFA92: Immediate Annuitants, females, amounts"
$ws.Range("E48").Value = "Mortality_tables/92series.xls"
$ws.Range("C49").Value = "IFA92: Immediate Annuitants, females, lives"
$ws.Range("E49").Value = "Mortality_tables/92series.xls"
$ws.Range("C50").Value = "IMA92: Immediate Annuitants, males, amounts"
$ws.Range("E50").Value = "Mortality_tables/92series.xls"
$ws.Range("C51").Value = "IMA92: Immediate Annuitants, males, lives"
$ws.Range("E51").Value = "Mortality_tables/92series.xls"
$ws.Range("C52").Value = "PFA92: Pensioners, females, amounts"
$ws.Range("E52").Value = "Mortality_tables/92series.xls"
$ws.Range("C53").Value = "PFL92: Pensioners, females, lives"
$ws.Range("E53").Value = "Mortality_tables/92series.xls"
$ws.Range("C54").Value = "PMA92: Pensioners, males, amounts"
$ws.Range("E54").Value = "Mortality_tables/92series.xls"
$ws.Range("C55").Value = "PML92: Pensioners, males, lives"
$ws.Range("E55").Value = "Mortality_tables/92series.xls"
$ws.Range("C56").Value = "RFV92: Retirement Annuitants, females, vested"
$ws.Range("E56").Value = "Mortality_tables/92series.xls"
$ws.Range("C57").Value = "RMV92: Retirement Annuitants, males, vested"
$ws.Range("E57").Value = "Mortality_tables/92series.xls"
$ws.Range("C58").Value = "TF92: Temporary Assurances, females, combined"
$ws.Range("E58").Value = "Mortality_tables/92series.xls"
$ws.Range("C59").Value = "TM92: Temporary Assurances, males, combined"
$ws.Range("E59").Value = "Mortality_tables/92series.xls"
$ws.Range("C60").Value = "WA92: Widows, amounts"
$ws.Range("E60").Value = "Mortality_tables/92series.xls"
$ws.Range("C61").Value = "WL92: Widows, lives"
$ws.Range("E61").Value = "Mortality_tables/92series.xls"

# Phase 3: fill Select Years (D) numeric column
$ws.Range("D46").Value = 2
$ws.Range("D47").Value = 2
$ws.Range("D48").Value = 1
$ws.Range("D49").Value = 1
$ws.Range("D50").Value = 1
$ws.Range("D51").Value = 1
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 0
$ws.Range("D55").Value = 0
$ws.Range("D56").Value = 0
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 5
$ws.Range("D59").Value = 5
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 0

# Update sheet view: zoom and selection to match the final saved view state
$ws.Activate()
$excel.ActiveWindow.Zoom = 87
$ws.Range("C66").Select()
